$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.023799657821655
$ws.Range("B1").Value = 0.8841099143028259
$ws.Range("C1").Value = 0.7697933912277222
$ws.Range("D1").Value = 4.195941925048828
$ws.Range("E1").Value = 1.674144506454468
